$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1380
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1380
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4140
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -4680

$ws.Range("H73").Value = 1380
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1380
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4140
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -6012

$ws.Range("H116").Value = 7998.375
$ws.Range("I116").Value = 1462
$ws.Range("K116").Value = 1462
$ws.Range("M116").Value = 1980

$ws.Range("H132").Value = 22224700
$ws.Range("I132").Value = 2711097.2
$ws.Range("J132").Value = 111120000
$ws.Range("K132").Value = 8133291.600000001
$ws.Range("L132").Value = 333360000
$ws.Range("M132").Value = -8130761.600000001
$ws.Range("N132").Value = -333365060

$ws.Range("H137").Value = 1792.4828
$ws.Range("I137").Value = 1345.5454
$ws.Range("J137").Value = 3197.1428
$ws.Range("K137").Value = 4036.6362
$ws.Range("L137").Value = 9591.428400000001
$ws.Range("M137").Value = -1486.6362
$ws.Range("N137").Value = -14691.4284

$ws.Range("H141").Value = 2098.75
$ws.Range("I141").Value = 1798.3334
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 5395.0002
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -215.0002000000004
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3070.1538
$ws.Range("I63").Value = 2216.6667
$ws.Range("J63").Value = 3801.7144
$ws.Range("K63").Value = 2216.6667
$ws.Range("L63").Value = 3801.7144
$ws.Range("M63").Value = -1530.6667
$ws.Range("N63").Value = -5173.7144

$ws.Range("H66").Value = 3070.1538
$ws.Range("I66").Value = 2216.6667
$ws.Range("J66").Value = 3801.7144
$ws.Range("K66").Value = 11083.3335
$ws.Range("L66").Value = 19008.572
$ws.Range("M66").Value = -7651.333500000001
$ws.Range("N66").Value = -25872.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3787.5881
$ws.Range("I134").Value = 2588.762
$ws.Range("J134").Value = 5724.154
$ws.Range("K134").Value = 7766.286
$ws.Range("L134").Value = 17172.462
$ws.Range("M134").Value = -5231.286
$ws.Range("N134").Value = -22242.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 904.75
$ws.Range("I16").Value = 731.2
$ws.Range("J16").Value = 1194
$ws.Range("K16").Value = 731.2
$ws.Range("L16").Value = 1194
$ws.Range("M16").Value = -444.2
$ws.Range("N16").Value = -1768

$ws.Range("H25").Value = 10544.286
$ws.Range("I25").Value = 1950
$ws.Range("J25").Value = 13982
$ws.Range("K25").Value = 1950
$ws.Range("L25").Value = 13982
$ws.Range("M25").Value = -1776
$ws.Range("N25").Value = -14330

$ws.Range("H31").Value = 30305130
$ws.Range("I31").Value = 83334300
$ws.Range("J31").Value = 2744.8096
$ws.Range("K31").Value = 83334300
$ws.Range("L31").Value = 2744.8096
$ws.Range("M31").Value = -83334005
$ws.Range("N31").Value = -3334.8096

$ws.Range("H34").Value = 30305130
$ws.Range("I34").Value = 83334300
$ws.Range("J34").Value = 2744.8096
$ws.Range("K34").Value = 83334300
$ws.Range("L34").Value = 2744.8096
$ws.Range("M34").Value = -83334098
$ws.Range("N34").Value = -3148.8096

$ws.Range("H107").Value = 1290.7354
$ws.Range("I107").Value = 455
$ws.Range("J107").Value = 4514.2856
$ws.Range("K107").Value = 455
$ws.Range("L107").Value = 4514.2856
$ws.Range("M107").Value = 1465
$ws.Range("N107").Value = -8354.285599999999

$ws.Range("H113").Value = 904.75
$ws.Range("I113").Value = 731.2
$ws.Range("J113").Value = 1194
$ws.Range("K113").Value = 731.2
$ws.Range("L113").Value = 1194
$ws.Range("M113").Value = 1438.8
$ws.Range("N113").Value = -5534

$ws.Range("H115").Value = 28500
$ws.Range("J115").Value = 28500
$ws.Range("L115").Value = 28500
$ws.Range("N115").Value = -30850

$ws.Range("H118").Value = 32063.092
$ws.Range("J118").Value = 32063.092
$ws.Range("L118").Value = 32063.092
$ws.Range("N118").Value = -35377.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 905.79
$ws.Range("J131").Value = 909.2708
$ws.Range("L131").Value = 2727.8124
$ws.Range("N131").Value = -12807.8124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1527.579
$ws.Range("I102").Value = 1193
$ws.Range("J102").Value = 2101.1428
$ws.Range("K102").Value = 1193
$ws.Range("L102").Value = 2101.1428
$ws.Range("M102").Value = 429
$ws.Range("N102").Value = -5345.1428

$ws.Range("H113").Value = 1484.1
$ws.Range("I113").Value = 1230.125
$ws.Range("K113").Value = 1230.125
$ws.Range("M113").Value = 939.875

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = ""

$ws.Range("H126").Value = 3560
$ws.Range("I126").Value = 4906
$ws.Range("J126").Value = 2214
$ws.Range("K126").Value = 14718
$ws.Range("L126").Value = 6642
$ws.Range("M126").Value = -12248
$ws.Range("N126").Value = -11582

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1006.06665
$ws.Range("I46").Value = 1420.2
$ws.Range("J46").Value = 799
$ws.Range("K46").Value = 1420.2
$ws.Range("L46").Value = 799
$ws.Range("M46").Value = -1232.2
$ws.Range("N46").Value = -1175

$ws.Range("H68").Value = 7721.316
$ws.Range("J68").Value = 2838.6924
$ws.Range("L68").Value = 2838.6924
$ws.Range("N68").Value = -4336.6924

$ws.Range("H71").Value = 7721.316
$ws.Range("J71").Value = 2838.6924
$ws.Range("L71").Value = 14193.462
$ws.Range("N71").Value = -21681.462

$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314

$ws.Range("H136").Value = 2544.1428
$ws.Range("I136").Value = 2126
$ws.Range("J136").Value = 3101.6667
$ws.Range("K136").Value = 6378
$ws.Range("L136").Value = 9305.000100000001
$ws.Range("M136").Value = -3828
$ws.Range("N136").Value = -14405.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

$ws.Range("H116").Value = 31000
$ws.Range("J116").Value = 31000
$ws.Range("L116").Value = 31000
$ws.Range("N116").Value = -40178

$ws.Range("H136").Value = 2820.4736
$ws.Range("I136").Value = 2845.1482
$ws.Range("J136").Value = 2759.9092
$ws.Range("K136").Value = 8535.444600000001
$ws.Range("L136").Value = 8279.7276
$ws.Range("M136").Value = -5985.444600000001
$ws.Range("N136").Value = -13379.7276
